$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for rows 2-62, replacing the previous
# "Strike#" based figures with real K counts (regenerated save_data).
$kValues = @{
    2 = 1
    3 = 0
    4 = 0
    5 = 1
    6 = 0
    7 = 4
    8 = 0
    9 = 2
    10 = 1
    11 = 3
    12 = 2
    13 = 3
    14 = 2
    15 = 3
    16 = 1
    17 = 1
    18 = 2
    19 = 3
    20 = 2
    21 = 3
    22 = 2
    23 = 3
    24 = 2
    25 = 2
    26 = 1
    27 = 3
    28 = 3
    29 = 2
    30 = 3
    31 = 3
    32 = 1
    33 = 3
    34 = 2
    35 = 3
    36 = 5
    37 = 2
    38 = 2
    39 = 3
    40 = 3
    41 = 1
    42 = 1
    43 = 4
    44 = 1
    45 = 3
    46 = 2
    47 = 2
    48 = 1
    49 = 0
    50 = 1
    51 = 3
    52 = 2
    53 = 3
    54 = 2
    55 = 2
    56 = 1
    57 = 2
    58 = 2
    59 = 3
    60 = 1
    61 = 2
    62 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
